$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 switches from a numeric date value to a literal text string "05-22-2024"
# Temporarily mark the cell as Text so Excel doesn't reinterpret the string as
# a date serial number, then restore the original date number format so the
# cell's style (s="2", numFmtId 180 dd/mm/yyyy) is left untouched.
$fmt = $ws.Range("C2").NumberFormat
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "05-22-2024"
$ws.Range("C2").NumberFormat = $fmt

# D2 value changes from 58 to 56
$ws.Range("D2").Value = 56
